# pilot 1: subject 1 + subject 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Print_Titles: $1:$3 -> $1:$1
# ---------------------------------------------------------------------
$nm = $wb.Names.Item("Elenco inventario!Print_Titles")
$nm.RefersTo = "='Elenco inventario'!`$1:`$1"

# ---------------------------------------------------------------------
# 2) Remove the two decorative blank rows at the top (old rows 1 & 2);
#    the old header row (3) becomes row 1 and the table moves up.
# ---------------------------------------------------------------------
$ws.Rows("1:2").Delete() | Out-Null

# ---------------------------------------------------------------------
# 3) Insert a new column before the old "TUTORIAL DURATION" column (G)
#    so the old column slides to H, ready to become "COMMENTS".
# ---------------------------------------------------------------------
$ws.Columns("G").Insert() | Out-Null

# ---------------------------------------------------------------------
# 4) Re-title the header cells.
# ---------------------------------------------------------------------
$ws.Range("B1").Value2 = "SUBJECT NUMBER"
$ws.Range("C1").Value2 = "AGE"
$ws.Range("D1").Value2 = "SEX"
$ws.Range("E1").Value2 = "LINK SENT"
$ws.Range("F1").Value2 = "EXP ERIMENT COMPLETED"
$ws.Range("G1").Value2 = "TUTORIAL DURATION (minutes)"
$ws.Range("H1").Value2 = "COMMENTS"

# ---------------------------------------------------------------------
# 5) Resize/refresh the table to match the new B1:H26 range + columns.
# ---------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B1:H26"))

# ---------------------------------------------------------------------
# 6) Fill in subject data rows.
# ---------------------------------------------------------------------
$ws.Range("B2").Value2 = 1
$ws.Range("C2").Value2 = 27
$ws.Range("D2").Value2 = "M"
$ws.Range("E2").Value2 = 0.42708333333333331
$ws.Range("F2").Value2 = 0.4458333333333333
$ws.Range("G2").Value2 = 11
$ws.Range("H2").Value2 = "it was not clear that the shuffles were starting from zero at every trial, so at the beginning he was trying to use as few as possible -> underline it more during the tutorial with feedback"

$ws.Range("B3").Value2 = 2
$ws.Range("C3").Value2 = 24
$ws.Range("D3").Value2 = "F"
$ws.Range("E3").Value2 = 0.4770833333333333
$ws.Range("F3").Value2 = 0.49861111111111112
$ws.Range("G3").Value2 = 13
$ws.Range("H3").Value2 = "bad quality images during tutorial (no feedback) + thought that using more shuffles resulted in more difficult following trials"

# ---------------------------------------------------------------------
# 7) Number formats.
# ---------------------------------------------------------------------
$ws.Range("E2:F3").NumberFormat = "h:mm;@"
$ws.Range("G2:G3").NumberFormat = "#,##0.00"

# ---------------------------------------------------------------------
# 8) Alignment / wrap for the comments column + row heights.
# ---------------------------------------------------------------------
$ws.Range("H2").HorizontalAlignment = -4131
$ws.Range("H2").VerticalAlignment = -4160
$ws.Range("H2").WrapText = $true

$ws.Range("H3").HorizontalAlignment = -4152
$ws.Range("H3").VerticalAlignment = -4108
$ws.Range("H3").WrapText = $true
$ws.Range("H3").IndentLevel = 1

$ws.Rows(2).RowHeight = 40.2
$ws.Rows(3).RowHeight = 45

# ---------------------------------------------------------------------
# 9) Column widths.
# ---------------------------------------------------------------------
$ws.Columns("H").ColumnWidth = 44.1796875

# ---------------------------------------------------------------------
# 10) Selection cosmetics (matches the authored file's cursor position).
# ---------------------------------------------------------------------
$ws.Range("B4").Select() | Out-Null
